# Lab5Rubric-CIS195.xlsx — "Updates to rubrics and instructions"
#
# 1. Rename Sheet1 -> Rubric, Sheet2 -> Score
# 2. Score sheet (Lab attempt): CSS width requirement wasn't met (C4: 5 -> 4),
#    add grader comments in column D for two rows that lost points.
# 3. Clean up the old bold+italic "Total" row style on both sheets so it
#    matches the plain bold style used elsewhere, and give sheet1 row 22 a
#    4th (italic, still-empty) cell to mirror the new column D on Score.
# 4. Flip which sheet is on top / selected, and tidy up the selections.

$wb = $excel.ActiveWorkbook

$wsRubric = $wb.Worksheets.Item(1)
$wsScore  = $wb.Worksheets.Item(2)

$wsRubric.Name = "Rubric"
$wsScore.Name  = "Score"

# --- Score sheet content updates ----------------------------------------

# "Set the body of the page to fixed width" actual score drops from 5 to 4
$wsScore.Range("C4").Value = 4

# Grader feedback notes for the two requirements that lost points
$wsScore.Range("D4").Value = "Should be in the external CSS"
$wsScore.Range("D5").Value = "Margin is fixed, but not width"

# --- Style cleanup: drop the bold+italic "Total" style in favor of bold ---

foreach ($ws in @($wsRubric, $wsScore)) {
    $ws.Range("A22:C22").Font.Italic = $false
}

# Sheet1 (Rubric) row 22 picks up an extra (empty, italic-styled) D cell so
# its "Total" row lines up with Score's new 4-column layout.
$wsRubric.Range("D22").Font.Italic = $true

# --- Page setup -----------------------------------------------------------
$wsRubric.PageSetup.Orientation = 1

# --- View / selection state -------------------------------------------

# Rubric: selection moves onto the Total row.
$wsRubric.Range("A22:C22").Select()

# Score becomes the active (front) sheet, zoomed to 120%, cursor on A14.
$wsScore.Activate()
$excel.ActiveWindow.Zoom = 120
$wsScore.Range("A14").Select()
